# change to 6 sats setting
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update satellite RAAN offsets (column D) for the new 6-satellite constellation setup
$ws.Range("D2").Value = -2
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = -6
$ws.Range("D5").Value = 6
$ws.Range("D6").Value = -10
$ws.Range("D7").Value = 10

# Font fell back from SimSun (unavailable on this machine) to Calibri.
# Only touch the populated cells so we don't materialize empty cells like I1.
$ws.Range("A1:H1").Font.Name = "Calibri"
$ws.Range("A2:I7").Font.Name = "Calibri"

# Move the selection, matching the saved state of the workbook.
$ws.Range("D8").Select()
